$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.590158771791153
$ws.Range("C2").Value = 0.6531155350378413
$ws.Range("D2").Value = 0.04807882700766442
$ws.Range("E2").Value = 0.1195447155408829
$ws.Range("F2").Value = 3.384970243692777
$ws.Range("J2").Value = 0.2274430023333451
$ws.Range("N2").Value = 1.80742071182388
$ws.Range("B3").Value = 2.437359848274127
$ws.Range("C3").Value = 0.6069069128337787
$ws.Range("D3").Value = 0.04793712971836328
$ws.Range("E3").Value = 0.1177158772567637
$ws.Range("F3").Value = 3.328819371221073
$ws.Range("J3").Value = 0.2223673216311681
$ws.Range("N3").Value = 1.822854897791139
$ws.Range("B4").Value = 2.345226233922006
$ws.Range("C4").Value = 0.5789541211847222
$ws.Range("D4").Value = 0.04786182353945279
$ws.Range("E4").Value = 0.1166596525626069
$ws.Range("F4").Value = 3.29650929750818
$ws.Range("J4").Value = 0.2193972202664654
$ws.Range("N4").Value = 1.832999409062658
$ws.Range("B5").Value = 2.308100492832921
$ws.Range("C5").Value = 0.5676671481161293
$ws.Range("D5").Value = 0.04783410073699557
$ws.Range("E5").Value = 0.1162459481883076
$ws.Range("F5").Value = 3.283883502477323
$ws.Range("J5").Value = 0.2182234284473736
$ws.Range("N5").Value = 1.837300489266397
$ws.Range("B6").Value = 2.301960999329992
$ws.Range("C6").Value = 0.565799194830447
$ws.Range("D6").Value = 0.04782967717819275
$ws.Range("E6").Value = 0.1161782609474393
$ws.Range("F6").Value = 3.281819546344749
$ws.Range("J6").Value = 0.2180307210996304
$ws.Range("N6").Value = 1.838024749266694
$ws.Range("B7").Value = 2.344723850566311
$ws.Range("C7").Value = 0.5788014818782017
$ws.Range("D7").Value = 0.04786143762296291
$ws.Range("E7").Value = 0.1166540055885825
$ws.Range("F7").Value = 3.296336837125978
$ws.Range("J7").Value = 0.2193812424570183
$ws.Range("N7").Value = 1.83305673939676
$ws.Range("B8").Value = 2.537121150666337
$ws.Range("C8").Value = 0.6370948311703728
$ws.Range("D8").Value = 0.04802755628597311
$ws.Range("E8").Value = 0.118900256927656
$ws.Range("F8").Value = 3.365157235577016
$ws.Range("J8").Value = 0.2256623641174258
$ws.Range("N8").Value = 1.812603304945874
$ws.Range("B9").Value = 2.928008787594422
$ws.Range("C9").Value = 0.7548098823408509
$ws.Range("D9").Value = 0.04844512876842089
$ws.Range("E9").Value = 0.1238371838021202
$ws.Range("F9").Value = 3.517501306237506
$ws.Range("J9").Value = 0.2391540164819759
$ws.Range("N9").Value = 1.777830493271864
$ws.Range("B10").Value = 3.223828490996596
$ws.Range("C10").Value = 0.8434814779429871
$ws.Range("D10").Value = 0.04880664150098823
$ws.Range("E10").Value = 0.1277933445842621
$ws.Range("F10").Value = 3.640317281162993
$ws.Range("J10").Value = 0.2498016674861674
$ws.Range("N10").Value = 1.75558719025139
$ws.Range("B11").Value = 3.360357790342334
$ws.Range("C11").Value = 0.8843206251108313
$ws.Range("D11").Value = 0.04898273104966222
$ws.Range("E11").Value = 0.129665624314562
$ws.Range("F11").Value = 3.698618511605162
$ws.Range("J11").Value = 0.254809598984167
$ws.Range("N11").Value = 1.74619687226803
$ws.Range("B12").Value = 3.412345509298802
$ws.Range("C12").Value = 0.8998594549823906
$ws.Range("D12").Value = 0.04905106291827011
$ws.Range("E12").Value = 0.130385128159233
$ws.Range("F12").Value = 3.721050406410285
$ws.Range("J12").Value = 0.256729930857432
$ws.Range("N12").Value = 1.742746698872011
$ws.Range("B13").Value = 3.401136172465272
$ws.Range("C13").Value = 0.8965095736118087
$ws.Range("D13").Value = 0.0490362733738543
$ws.Range("E13").Value = 0.1302297014287035
$ws.Range("F13").Value = 3.716203452947582
$ws.Range("J13").Value = 0.256315283317889
$ws.Range("N13").Value = 1.743485036492572
$ws.Range("B14").Value = 3.364629067966803
$ws.Range("C14").Value = 0.885597522160765
$ws.Range("D14").Value = 0.04898831980164786
$ws.Range("E14").Value = 0.129724607228269
$ws.Range("F14").Value = 3.700456864757228
$ws.Range("J14").Value = 0.2549671042745274
$ws.Range("N14").Value = 1.74591089946135
$ws.Range("B15").Value = 3.342304951657525
$ws.Range("C15").Value = 0.8789232560890241
$ws.Range("D15").Value = 0.0489591611428537
$ws.Range("E15").Value = 0.1294165936187994
$ws.Range("F15").Value = 3.690857933777068
$ws.Range("J15").Value = 0.2541444324838125
$ws.Range("N15").Value = 1.74741061086462
$ws.Range("B16").Value = 3.214945856720988
$ws.Range("C16").Value = 0.8408227932679324
$ws.Range("D16").Value = 0.04879536583211674
$ws.Range("E16").Value = 0.1276724537439549
$ws.Range("F16").Value = 3.636556498965888
$ws.Range("J16").Value = 0.2494777196240534
$ws.Range("N16").Value = 1.756215610390385
$ws.Range("B17").Value = 3.137320745216755
$ws.Range("C17").Value = 0.8175791947598441
$ws.Range("D17").Value = 0.04869784666498589
$ws.Range("E17").Value = 0.1266211260672314
$ws.Range("F17").Value = 3.603870524313891
$ws.Range("J17").Value = 0.2466571365859238
$ws.Range("N17").Value = 1.761804434660164
$ws.Range("B18").Value = 3.092856942101662
$ws.Range("C18").Value = 0.8042572348391559
$ws.Range("D18").Value = 0.04864285274443958
$ws.Range("E18").Value = 0.1260232594106547
$ws.Range("F18").Value = 3.585298998594453
$ws.Range("J18").Value = 0.2450502535665038
$ws.Range("N18").Value = 1.765087465476512
$ws.Range("B19").Value = 3.077833724529171
$ws.Range("C19").Value = 0.7997546924992207
$ws.Range("D19").Value = 0.04862442177156012
$ws.Range("E19").Value = 0.1258220028051191
$ws.Range("F19").Value = 3.579050113020287
$ws.Range("J19").Value = 0.2445088331996459
$ws.Range("N19").Value = 1.766210775429826
$ws.Range("B20").Value = 3.145564978607922
$ws.Range("C20").Value = 0.8200486238512212
$ws.Range("D20").Value = 0.04870811442434686
$ws.Range("E20").Value = 0.1267323343816571
$ws.Range("F20").Value = 3.607326310109187
$ws.Range("J20").Value = 0.2469557920352941
$ws.Range("N20").Value = 1.761202399370745
$ws.Range("B21").Value = 3.375344256750054
$ws.Range("C21").Value = 0.888800635602081
$ws.Range("D21").Value = 0.04900236031050298
$ws.Range("E21").Value = 0.1298726797554508
$ws.Range("F21").Value = 3.70507235992622
$ws.Range("J21").Value = 0.2553624449434437
$ws.Range("N21").Value = 1.74519548690229
$ws.Range("B22").Value = 3.527193906445973
$ws.Range("C22").Value = 0.9341658166572415
$ws.Range("D22").Value = 0.04920427762484181
$ws.Range("E22").Value = 0.131986385072878
$ws.Range("F22").Value = 3.771023564472614
$ws.Range("J22").Value = 0.2609963531016888
$ws.Range("N22").Value = 1.735350837760521
$ws.Range("B23").Value = 3.445993796063135
$ws.Range("C23").Value = 0.9099134716586263
$ws.Range("D23").Value = 0.04909563855219545
$ws.Range("E23").Value = 0.130852626609709
$ws.Range("F23").Value = 3.735633279034403
$ws.Range("J23").Value = 0.2579765436852881
$ws.Range("N23").Value = 1.740548336851774
$ws.Range("B24").Value = 3.14183725055284
$ws.Range("C24").Value = 0.8189320670900884
$ws.Range("D24").Value = 0.04870346903219769
$ws.Range("E24").Value = 0.1266820366888517
$ws.Range("F24").Value = 3.605763264225288
$ws.Range("J24").Value = 0.2468207241189049
$ws.Range("N24").Value = 1.761474361750018
$ws.Range("B25").Value = 2.820772391448827
$ws.Range("C25").Value = 0.7225891535520645
$ws.Range("D25").Value = 0.0483224656412311
$ws.Range("E25").Value = 0.1224441832956735
$ws.Range("F25").Value = 3.474397201050635
$ws.Range("J25").Value = 0.2353764591542671
$ws.Range("N25").Value = 1.786660944662046
